$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# --- Update the existing data row (row 2): the Start/End dates change
# (14/04/2021 -> 14/06/2021, 20/05/2021 -> 20/08/2021) for every row that
# references them.
$ws.Range("H2").Value = "'14/06/2021"
$ws.Range("I2").Value = "'20/08/2021"

# --- Add a new row 4 first so the "Cypress" shared string is registered
# before the longer description string (keeps shared-string table order
# the same as Excel would produce when typing the new title first).
$ws.Range("A4:P4").Value = $ws.Range("A2:P2").Value()
$ws.Range("A4").Value = "Cypress"

# --- Add a new row 3: a copy of row 2, but offering Cypress training
# (still listed under the "Selenium" title, per the source data).
$ws.Range("A3:P3").Value = $ws.Range("A2:P2").Value()
$ws.Range("B3").Value = "Would like to provide Cypress training for beginners"
$ws.Range("B4").Value = "Would like to provide Cypress training for beginners"

# --- Re-apply the per-column formatting that plain value assignment does not
# carry over: wrap text on the Description column, and a text/quote-prefix
# style on the date & time columns (matches the style used on row 2).
foreach ($r in 3, 4) {
    $ws.Range("B$r").WrapText = $true

    $ws.Range("H$r").Value = "'" + $ws.Range("H$r").Value()
    $ws.Range("I$r").Value = "'" + $ws.Range("I$r").Value()
    $ws.Range("K$r").Value = "'" + $ws.Range("K$r").Value()
    $ws.Range("L$r").Value = "'" + $ws.Range("L$r").Value()

    $ws.Rows.Item($r).RowHeight = 43.2
}

# --- Update the sheet view: scroll back to column A and select B4.
$ws.Activate()
$ws.Range("B4").Select()
